# Updates cryptocurrency price/volume data in the worksheet to reflect the
# latest scrape, and swaps the Toncoin / InjectiveProtocol rows (33 and 34)
# since their ranking order changed.
#
# Leading single-quote on each assigned string forces Excel to treat the
# value as literal text (quote-prefix), preventing values like "384.18" or
# "51.602.97" from being auto-converted into numbers/dates, and keeping the
# cell's number format as "General" (matching the original workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'51.602.97"
$ws.Range("E2").Value = "'  +1.01%  "

# Row 3
$ws.Range("D3").Value = "'3.038.20"
$ws.Range("E3").Value = "'  +2.59%  "

# Row 4
$ws.Range("E4").Value = "'  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'384.18"
$ws.Range("E5").Value = "'  +1.09%  "

# Row 6
$ws.Range("D6").Value = "'102.58"
$ws.Range("E6").Value = "'  +0.52%  "

# Row 7
$ws.Range("D7").Value = "'0.543"
$ws.Range("E7").Value = "'  -0.22%  "

# Row 8
$ws.Range("E8").Value = "'  -0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = "'  -0.34%  "

# Row 10
$ws.Range("D10").Value = "'36.79"
$ws.Range("E10").Value = "'  +0.88%  "

# Row 11
$ws.Range("E11").Value = "'  +0.18%  "

# Row 12
$ws.Range("E12").Value = "'  +0.94%  "

# Row 13
$ws.Range("D13").Value = "'3.516.22"
$ws.Range("E13").Value = "'  +2.20%  "

# Row 14
$ws.Range("D14").Value = "'18.71"
$ws.Range("E14").Value = "'  +2.07%  "

# Row 15
$ws.Range("D15").Value = "'7.77"
$ws.Range("E15").Value = "'  -0.57%  "

# Row 16
$ws.Range("D16").Value = "'3.041.94"
$ws.Range("E16").Value = "'  +3.12%  "

# Row 17
$ws.Range("D17").Value = "'0.972"
$ws.Range("E17").Value = "'  -2.64%  "

# Row 18
$ws.Range("D18").Value = "'10.57"
$ws.Range("E18").Value = "'  -10.91%  "

# Row 19
$ws.Range("D19").Value = "'51.636.91"
$ws.Range("E19").Value = "'  +0.87%  "

# Row 20
$ws.Range("E20").Value = "'  -0.59%  "

# Row 21
$ws.Range("D21").Value = "'12.37"
$ws.Range("E21").Value = "'  -0.20%  "

# Row 22
$ws.Range("E22").Value = "'  +0.40%  "

# Row 23
$ws.Range("D23").Value = "'69.96"
$ws.Range("E23").Value = "'  -0.12%  "

# Row 24
$ws.Range("D24").Value = "'267.03"
$ws.Range("E24").Value = "'  +0.01%  "

# Row 25
$ws.Range("D25").Value = "'3.16"
$ws.Range("E25").Value = "'  -3.49%  "

# Row 26
$ws.Range("D26").Value = "'8.40"
$ws.Range("E26").Value = "'  +6.20%  "

# Row 27
$ws.Range("D27").Value = "'7.41"
$ws.Range("E27").Value = "'  +4.04%  "

# Row 28
$ws.Range("E28").Value = "'  +3.43%  "

# Row 29
$ws.Range("D29").Value = "'26.40"
$ws.Range("E29").Value = "'  +2.08%  "

# Row 30
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "'  -0.07%  "

# Row 31
$ws.Range("E31").Value = "'  -2.70%  "

# Row 32
$ws.Range("E32").Value = "'  -0.19%  "

# Row 33
$ws.Range("B33").Value = "'Toncoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'2.07"
$ws.Range("E33").Value = "'  +0.45%  "

# Row 34
$ws.Range("B34").Value = "'InjectiveProtocol"
$ws.Range("C34").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'34.03"
$ws.Range("E34").Value = "'  -1.08%  "

# Row 35
$ws.Range("D35").Value = "'50.62"
$ws.Range("E35").Value = "'  -1.21%  "

# Row 36
$ws.Range("D36").Value = "'0.0445"
$ws.Range("E36").Value = "'  +2.25%  "

# Row 37
$ws.Range("E37").Value = "'  -0.12%  "

# Row 38
$ws.Range("E38").Value = "'  +3.52%  "

# Row 39
$ws.Range("D39").Value = "'0.286"
$ws.Range("E39").Value = "'  +5.38%  "

# Row 40
$ws.Range("D40").Value = "'16.97"
$ws.Range("E40").Value = "'  +2.84%  "

# Row 41
$ws.Range("E41").Value = "'  +1.70%  "

# Row 42
$ws.Range("D42").Value = "'127.81"
$ws.Range("E42").Value = "'  +1.99%  "

# Row 43
$ws.Range("E43").Value = "'  -0.50%  "

# Row 44
$ws.Range("D44").Value = "'2.52"
$ws.Range("E44").Value = "'  +0.68%  "

# Row 45
$ws.Range("E45").Value = "'  +3.21%  "

# Row 46
$ws.Range("D46").Value = "'21.67"
$ws.Range("E46").Value = "'  +0.72%  "

# Row 47
$ws.Range("E47").Value = "'  +3.95%  "

# Row 48
$ws.Range("D48").Value = "'2.07"
$ws.Range("E48").Value = "'  +2.41%  "

# Row 49
$ws.Range("D49").Value = "'2.032.81"
$ws.Range("E49").Value = "'  -0.87%  "

# Row 50
$ws.Range("D50").Value = "'3.341.51"
$ws.Range("E50").Value = "'  +2.54%  "

# Row 51
$ws.Range("D51").Value = "'0.207"
$ws.Range("E51").Value = "'  +7.06%  "
